$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B89").Value = 0
$ws.Range("C89").Value = "key"
